# Mark every "realizacja" (completion) checkbox in the first scoring table
# (XML->HTML, rows 3-15) as fully done (1 instead of 0). The dependent
# "pkt" column D (B*C) and the summary totals (D16, D29, D37) recalc
# automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2 projekt")

$ws.Range("C3:C15").Value = 1

# Restore the cursor/selection to the cell the author last had selected.
$ws.Range("E20").Select() | Out-Null
